# Update the "Instructions" sheet of the FWC Calculator data template:
#  - Insert a new requirement bullet stating the data must cover exactly one storm event.
#  - Clarify that flow/pollutant measurements must be >= 0 (not strictly > 0).

$wb = $excel.ActiveWorkbook
$instructions = $wb.Worksheets.Item("Instructions")

# Insert a new row above the existing "Must contain exactly three sheets..." bullet
# (currently row 15) so the new bullet becomes the first item under "Data Requirements".
$instructions.Rows.Item(15).Insert()
$instructions.Range("A15").Value = "  * Must contain data for exactly one storm event."

# The "greater than zero" bullet (now shifted down to row 31) should read
# "greater than or equal to zero" instead.
$instructions.Range("A31").Value = "  * All flow rate and pollutant measurements must be greater than or equal to zero."

# Restore the active selection on the Instructions sheet.
$instructions.Activate()
$instructions.Range("A12").Select()
